$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Copy the header style from an existing header cell (D1) to the new headers
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in Execution Time (ms) values
$ws.Range("E2").Value = 7.384800002910197
$ws.Range("E3").Value = 6.279799999902025
$ws.Range("E4").Value = 25.28699999675155
$ws.Range("E5").Value = 4.357500001788139
$ws.Range("E6").Value = 1.843600010033697

# Fill in Memory Usage (B) values
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
